$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K3:K6").FormulaArray = "=LAMBDA(OriginalText,ReplacementMap, LET(OriginalText2, B3:B6, ReplacementMap2, D3:E8, FirstRow, INDEX(OriginalText, 1, 1), Seq, SEQUENCE(ROWS(ReplacementMap)), Result, BYROW(OriginalText, LAMBDA(CurrentOriginalText, REDUCE(CurrentOriginalText, Seq, LAMBDA(Acc,Curr, SUBSTITUTE(Acc, INDEX(ReplacementMap, Curr, 1), INDEX(ReplacementMap, Curr, 2)))))), Result))(B3:B6, D3:E8)"
